$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Column widths: A and B split apart and widened, new column D added.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 54.33
$ws.Columns.Item(2).ColumnWidth = 57.83
$ws.Columns.Item(4).ColumnWidth = 78.666666666666671

# ---------------------------------------------------------------------------
# New MOR-table block: rows 8-10, merged A8:A10 / C8:C10, new "Doing" status
# and a new wide note column (D). Formats are copied from the existing
# COW-table block (rows 5-7) before row 6's own format is tweaked below.
# ---------------------------------------------------------------------------
$ws.Range("A5").Copy()
$ws.Range("A8:A10").PasteSpecial(-4122)
$ws.Range("C5").Copy()
$ws.Range("C8:C10").PasteSpecial(-4122)
$ws.Range("B6").Copy()
$ws.Range("B8").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Row 6 is now shorter (the COW-table "Observation" row) and vertically
# centred.
# ---------------------------------------------------------------------------
$ws.Rows.Item(6).RowHeight = 80
$ws.Range("B6").VerticalAlignment = -4108

$ws.Range("C8").Value = "Doing"
$ws.Range("A8").Value = "MOR Table (Applicable for Format-version : 2 only) "
$ws.Range("B8").Value = " Update a record will create two new files, a file which contains file path of the previous file and the position of the record in the table which is to be deleted (or deleted for updating the table) , and the second file which only contains the updated record."

$d8 = $ws.Range("D8")
$d8.Value = "The format-version property in an Apache Iceberg table specifies the version of the Iceberg table format that the table uses. The current version of the Iceberg table format is 2. Version 1 of the Iceberg table format does not support row-level deletes.  Version 2 of the Iceberg table format adds support for row-level deletes. This means that you can delete or replace individual rows in immutable data files without rewriting the files."
$d8.HorizontalAlignment = -4108
$d8.VerticalAlignment = -4108
$d8.WrapText = $true
$d8.ShrinkToFit = $true
$d8.Font.Size = 12
$d8.Font.Name = "Helvetica Neue"

$ws.Rows.Item(8).RowHeight = 100

$ws.Range("A8:A10").Merge()
$ws.Range("C8:C10").Merge()

# ---------------------------------------------------------------------------
# View tweaks: zoom level and active selection moved onto the new column.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 104
$ws.Range("D11").Select()
